$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text (shared-string) values first, in the same order as the original edit,
# so the shared strings table ends up in the same order.
$ws.Range("F27").Value = "Exp 29.png"
$ws.Range("A27").Value = "Exp 29"
$ws.Range("A28").Value = "Exp 30"
$ws.Range("F28").Value = "Exp 30.png"
$ws.Range("A29").Value = "Exp 31"
$ws.Range("F29").Value = "Exp 31.png"

$ws.Range("D27").Value = "Local"
$ws.Range("D28").Value = "Local"
$ws.Range("D29").Value = "Local"

# Row 27 - Exp 29
$ws.Range("B27").Value = 0.1
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = -1
$ws.Range("G27").Value = 66.32
$ws.Range("H27").Value = 67.15
$ws.Range("I27").Value = 54.64
$ws.Range("J27").Value = 48.68
$ws.Range("K27").Value = 45.74

# Row 28 - Exp 30
$ws.Range("B28").Value = 0.3
$ws.Range("C28").Value = 30
$ws.Range("E28").Value = -1
$ws.Range("G28").Value = 63.38
$ws.Range("H28").Value = 62
$ws.Range("I28").Value = 79.57
$ws.Range("J28").Value = 45.89
$ws.Range("K28").Value = 48.16

# Row 29 - Exp 31
$ws.Range("B29").Value = 0.1
$ws.Range("C29").Value = 30
$ws.Range("E29").Value = -1
$ws.Range("G29").Value = 22.19
$ws.Range("H29").Value = 17.31
$ws.Range("I29").Value = 53.17
$ws.Range("J29").Value = 9.26
$ws.Range("K29").Value = 12.43

# Match the centered styling used by the rest of the table (A:E and G:K columns)
$ws.Range("A27:E29").HorizontalAlignment = -4108
$ws.Range("G27:K29").HorizontalAlignment = -4108

# Update view: scroll position and selection to mirror the diff
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("K30").Select()
